$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 2,3,4,5

foreach ($r in $rows) {
    $ws.Range("C$r").Value() = "Unknown Title"
    $ws.Range("D$r").Value() = "Unknown Abstract"
    $ws.Range("E$r").Value() = "[]"
    $ws.Range("F$r").Value() = "not found"
    $ws.Range("G$r").Value() = "N/A"

    $ws.Range("H$r").NumberFormat = "@"
    $ws.Range("H$r").Value() = "1970-01-01"

    $ws.Range("I$r").Value() = ""
    $ws.Range("J$r").Value() = ""
}
